# MuS_Data.xlsx – "Improvements gemengt und data collected."
# Adds a fourth (Z:AC) mini-table of collected data to rows 1-4 of
# Tabelle1, alongside the existing J:N / P:T tables, and updates the
# sheet selection to reflect where the user was working (T20:U20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block in columns Z:AC (rows 1-4): discipline label + three
# measured values, mirroring the layout already used by the P:T table.
$ws.Range("Z1").Value = "LPT"
$ws.Range("AA1").Value = 2370
$ws.Range("AB1").Value = 2788
$ws.Range("AC1").Value = 5158

$ws.Range("Z2").Value = "SPT"
$ws.Range("AA2").Value = 1873
$ws.Range("AB2").Value = 2788
$ws.Range("AC2").Value = 4661

$ws.Range("Z3").Value = "LIFO"
$ws.Range("AA3").Value = 2148
$ws.Range("AB3").Value = 2788
$ws.Range("AC3").Value = 4937

$ws.Range("Z4").Value = "FIFO"
$ws.Range("AA4").Value = 2075
$ws.Range("AB4").Value = 2788
$ws.Range("AC4").Value = 4864

# Leave the selection where the author last worked before saving.
$ws.Range("T20:U20").Select()
